# Append plot-center survey rows captured at Loon Lake to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "plot"
$ws.Range("B1").Value = "x"
$ws.Range("C1").Value = "y"
$ws.Range("A1:C1").Style = "Normal"

# Plot 1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = -122.58965
$ws.Range("C2").Value = 49.29969
$ws.Range("A2:C2").Style = "Normal"

# Plot 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = -122.58941
$ws.Range("C3").Value = 49.30012

# Plot 3
$ws.Range("A4").Value = 3
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = -122.58847
$ws.Range("C4").Value = 49.30004

# Plot 4
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = -122.58857
$ws.Range("C5").Value = 49.29956

# Plot 5
$ws.Range("A6").Value = 5
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = -122.58965
$ws.Range("C6").Value = 49.29908

# Plot 6
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = -122.59019
$ws.Range("C7").Value = 49.29994

$ws.Range("A1").Select()
